$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New participants data (rows 17-20)
# columns: A=Nb, B=Name, C=Phone numbers, D=Amount paid, E=Method, F=Amount reimbursed, G=Amount left to reimburse
$rows = @(
    @{ Nb = 16; Name = "Ndjiki";         Phone = 690839895; Paid = 300; Method = "cash"; Reimbursed = 0;   HasFormula = $true  },
    @{ Nb = 17; Name = "Dasse";          Phone = 694993298; Paid = 300; Method = "cash"; Reimbursed = 0;   HasFormula = $true  },
    @{ Nb = 18; Name = "Gaetan Aymar";   Phone = 691674935; Paid = 500; Method = "cash"; Reimbursed = 100; HasFormula = $true  },
    @{ Nb = 19; Name = "Amombo Ngongo";  Phone = 691569975; Paid = 500; Method = "cash"; Reimbursed = 200; HasFormula = $false }
)

$startRow = 17
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Nb
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 3).Value = $row.Phone
    $ws.Cells.Item($r, 4).Value = $row.Paid
    $ws.Cells.Item($r, 5).Value = $row.Method
    $ws.Cells.Item($r, 6).Value = $row.Reimbursed
    if ($row.HasFormula) {
        $ws.Cells.Item($r, 7).Formula = "=D$r-F$r-300"
    }
    $r = $r + 1
}

# Update the view: scroll so row 4 is at top, and select C20 (last entered cell's column)
$ws.Activate()
$ws.Range("C20").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
